# Add season-record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the three new columns, placed right after the existing
# "Unnamed: 28" column (AC). Copy the formatting of AC1 (bold/border/
# centered header style) onto AD1:AF1 before writing the text so the
# new header cells reuse the same style as the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Season record values for every player row (2-56): 84 wins, 78 losses,
# 0 ties for this team's season.
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 30).Value = 84
    $ws.Cells.Item($r, 31).Value = 78
    $ws.Cells.Item($r, 32).Value = 0
}
